$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 962
$ws.Range("I4").Value = 962
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 962
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -848
$ws.Range("N4").ClearContents()
$ws.Range("H11").Value = 90917.82000000001
$ws.Range("I11").Value = 90917.82000000001
$ws.Range("K11").Value = 90917.82000000001
$ws.Range("M11").Value = -90777.82000000001
$ws.Range("H32").Value = 334.70587
$ws.Range("I32").Value = 330.2
$ws.Range("J32").Value = 336.58334
$ws.Range("K32").Value = 330.2
$ws.Range("L32").Value = 336.58334
$ws.Range("M32").Value = -4.199999999999989
$ws.Range("N32").Value = -988.58334
$ws.Range("H41").Value = 709.3461
$ws.Range("I41").Value = 957.4286
$ws.Range("J41").Value = 617.9474
$ws.Range("K41").Value = 957.4286
$ws.Range("L41").Value = 617.9474
$ws.Range("M41").Value = -517.4286
$ws.Range("N41").Value = -1497.9474
$ws.Range("H53").Value = 484.4
$ws.Range("I53").Value = 500.81818
$ws.Range("J53").Value = 471.5
$ws.Range("K53").Value = 500.81818
$ws.Range("L53").Value = 471.5
$ws.Range("M53").Value = 136.18182
$ws.Range("N53").Value = -1745.5
$ws.Range("H111").Value = 12799.909
$ws.Range("I111").Value = 13581.9
$ws.Range("J111").Value = 4980
$ws.Range("K111").Value = 40745.7
$ws.Range("L111").Value = 14940
$ws.Range("M111").Value = -37678.7
$ws.Range("N111").Value = -21074
$ws.Range("H135").Value = 1384.3948
$ws.Range("I135").Value = 917.6923
$ws.Range("J135").Value = 1627.08
$ws.Range("K135").Value = 8259.2307
$ws.Range("L135").Value = 14643.72
$ws.Range("M135").Value = -5724.2307
$ws.Range("N135").Value = -19713.72
$ws.Range("H137").Value = 1496.9286
$ws.Range("I137").Value = 1644.2941
$ws.Range("J137").Value = 1269.1818
$ws.Range("K137").Value = 4932.8823
$ws.Range("L137").Value = 3807.5454
$ws.Range("M137").Value = -2382.8823
$ws.Range("N137").Value = -8907.545399999999
$ws.Range("H138").Value = 4693.381
$ws.Range("I138").Value = 2877.1
$ws.Range("J138").Value = 6344.5454
$ws.Range("K138").Value = 8631.299999999999
$ws.Range("L138").Value = 19033.6362
$ws.Range("M138").Value = -3491.299999999999
$ws.Range("N138").Value = -29313.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 125250940
$ws.Range("I110").Value = 125250940
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 125250940
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -125248895
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 73841.42999999999
$ws.Range("I20").Value = 102278
$ws.Range("J20").Value = 2750
$ws.Range("K20").Value = 102278
$ws.Range("L20").Value = 2750
$ws.Range("M20").Value = -102031
$ws.Range("N20").Value = -3244
$ws.Range("H64").Value = 636.6
$ws.Range("I64").Value = 453.2
$ws.Range("J64").Value = 820
$ws.Range("K64").Value = 453.2
$ws.Range("L64").Value = 820
$ws.Range("M64").Value = -228.2
$ws.Range("N64").Value = -1270
$ws.Range("H67").Value = 636.6
$ws.Range("I67").Value = 453.2
$ws.Range("J67").Value = 820
$ws.Range("K67").Value = 453.2
$ws.Range("L67").Value = 820
$ws.Range("M67").Value = 326.8
$ws.Range("N67").Value = -2380
$ws.Range("H86").Value = 55574.895
$ws.Range("I86").Value = 115401.22
$ws.Range("J86").Value = 1731.2
$ws.Range("K86").Value = 115401.22
$ws.Range("L86").Value = 1731.2
$ws.Range("M86").Value = -114278.22
$ws.Range("N86").Value = -3977.2
$ws.Range("H89").Value = 55574.895
$ws.Range("I89").Value = 115401.22
$ws.Range("J89").Value = 1731.2
$ws.Range("K89").Value = 577006.1
$ws.Range("L89").Value = 8656
$ws.Range("M89").Value = -571390.1
$ws.Range("N89").Value = -19888

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 14999
$ws.Range("I45").Value = 14998
$ws.Range("K45").Value = 14998
$ws.Range("M45").Value = -14405
$ws.Range("H132").Value = 166675580
$ws.Range("I132").Value = 200012050
$ws.Range("J132").Value = 125005000
$ws.Range("K132").Value = 600036150
$ws.Range("L132").Value = 375015000
$ws.Range("M132").Value = -600033620
$ws.Range("N132").Value = -375020060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1379.5333
$ws.Range("I5").Value = 1222.5714
$ws.Range("J5").Value = 1516.875
$ws.Range("K5").Value = 3667.7142
$ws.Range("L5").Value = 4550.625
$ws.Range("M5").Value = -3555.7142
$ws.Range("N5").Value = -4774.625
$ws.Range("H40").Value = 358.25
$ws.Range("I40").Value = 123.42857
$ws.Range("J40").Value = 2002
$ws.Range("K40").Value = 493.71428
$ws.Range("L40").Value = 8008
$ws.Range("M40").Value = -424.71428
$ws.Range("N40").Value = -8146
$ws.Range("H131").Value = 827.62244
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 837.96844
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 2513.90532
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -12593.90532
$ws.Range("H135").Value = 1379.5333
$ws.Range("I135").Value = 1222.5714
$ws.Range("J135").Value = 1516.875
$ws.Range("K135").Value = 11003.1426
$ws.Range("L135").Value = 13651.875
$ws.Range("M135").Value = -8468.142600000001
$ws.Range("N135").Value = -18721.875
$ws.Range("H140").Value = 4890.933
$ws.Range("I140").Value = 6452.3687
$ws.Range("K140").Value = 19357.1061
$ws.Range("M140").Value = -14177.1061

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 71136.2
$ws.Range("I70").Value = 99549.57000000001
$ws.Range("J70").Value = 4838.3335
$ws.Range("K70").Value = 99549.57000000001
$ws.Range("L70").Value = 4838.3335
$ws.Range("M70").Value = -99279.57000000001
$ws.Range("N70").Value = -5378.3335
$ws.Range("H73").Value = 71136.2
$ws.Range("I73").Value = 99549.57000000001
$ws.Range("J73").Value = 4838.3335
$ws.Range("K73").Value = 99549.57000000001
$ws.Range("L73").Value = 4838.3335
$ws.Range("M73").Value = -98613.57000000001
$ws.Range("N73").Value = -6710.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H46").Value = 6121.3335
$ws.Range("I46").Value = 1479.2858
$ws.Range("J46").Value = 12620.2
$ws.Range("K46").Value = 1479.2858
$ws.Range("L46").Value = 12620.2
$ws.Range("M46").Value = -1291.2858
$ws.Range("N46").Value = -12996.2
$ws.Range("H55").Value = 291974.88
$ws.Range("I55").Value = 494822.97
$ws.Range("J55").Value = 380.75
$ws.Range("K55").Value = 494822.97
$ws.Range("L55").Value = 380.75
$ws.Range("M55").Value = -494649.97
$ws.Range("N55").Value = -726.75
$ws.Range("H61").Value = 1681.0435
$ws.Range("I61").Value = 1526.9333
$ws.Range("J61").Value = 1970
$ws.Range("K61").Value = 1526.9333
$ws.Range("L61").Value = 1970
$ws.Range("M61").Value = -1324.9333
$ws.Range("N61").Value = -2374
$ws.Range("H113").Value = 1681.0435
$ws.Range("I113").Value = 1526.9333
$ws.Range("J113").Value = 1970
$ws.Range("K113").Value = 1526.9333
$ws.Range("L113").Value = 1970
$ws.Range("M113").Value = 643.0667000000001
$ws.Range("N113").Value = -6310

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 7325
$ws.Range("I20").Value = 2000
$ws.Range("J20").Value = 9100
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 9100
$ws.Range("M20").Value = -1760
$ws.Range("N20").Value = -9580
$ws.Range("H54").Value = 6928.9165
$ws.Range("J54").Value = 6916.091
$ws.Range("L54").Value = 6916.091
$ws.Range("N54").Value = -7956.091
$ws.Range("H113").Value = 712.9375
$ws.Range("J113").Value = 1100.6
$ws.Range("L113").Value = 3301.8
$ws.Range("N113").Value = -7641.799999999999
$ws.Range("H122").Value = 1276.48
$ws.Range("I122").Value = 1221.5
$ws.Range("J122").Value = 1346.4546
$ws.Range("K122").Value = 3664.5
$ws.Range("L122").Value = 4039.3638
$ws.Range("M122").Value = -1214.5
$ws.Range("N122").Value = -8939.363799999999
$ws.Range("H132").Value = 2654.8276
$ws.Range("I132").Value = 2621.7036
$ws.Range("J132").Value = 3102
$ws.Range("K132").Value = 7865.110799999999
$ws.Range("L132").Value = 9306
$ws.Range("M132").Value = -5335.110799999999
$ws.Range("N132").Value = -14366
